$wb = $excel.ActiveWorkbook

# --- Update GbPbT sheet: replace NOx (row 5) and OC (row 10) formulas with literal 0 ---
$ws = $wb.Worksheets.Item("GbPbT")
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0

# --- Add note to About sheet ---
$about = $wb.Worksheets.Item("About")
$about.Range("A15").Value = "We have customized the India EPS to use values of 0 for Nox and OC."
